$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 65.2
$ws.Range("J9").Value = 88
$ws.Range("L9").Value = 88
$ws.Range("N9").Value = -426
$ws.Range("H15").Value = 2391
$ws.Range("I15").Value = 2391
$ws.Range("K15").Value = 7173
$ws.Range("M15").Value = -7004
$ws.Range("H17").Value = 1071.1818
$ws.Range("J17").Value = 1076.8462
$ws.Range("L17").Value = 3230.5386
$ws.Range("N17").Value = -3566.5386
$ws.Range("H18").Value = 1296.3334
$ws.Range("I18").Value = 1296.3334
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1296.3334
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1012.3334
$ws.Range("N18").ClearContents()
$ws.Range("H43").Value = 4200.091
$ws.Range("J43").Value = 4133.6665
$ws.Range("L43").Value = 4133.6665
$ws.Range("N43").Value = -4271.6665
$ws.Range("H51").Value = 2994.8455
$ws.Range("I51").Value = 2994.318
$ws.Range("K51").Value = 2994.318
$ws.Range("M51").Value = -2510.318
$ws.Range("H70").Value = 4721.3887
$ws.Range("J70").Value = 4811.5625
$ws.Range("L70").Value = 14434.6875
$ws.Range("N70").Value = -14974.6875
$ws.Range("H73").Value = 4721.3887
$ws.Range("J73").Value = 4811.5625
$ws.Range("L73").Value = 14434.6875
$ws.Range("N73").Value = -16306.6875
$ws.Range("H86").Value = 8394.643
$ws.Range("I86").Value = 9010.875
$ws.Range("J86").Value = 7573
$ws.Range("K86").Value = 9010.875
$ws.Range("L86").Value = 7573
$ws.Range("M86").Value = -7887.875
$ws.Range("N86").Value = -9819
$ws.Range("H88").Value = 3800.3333
$ws.Range("J88").Value = 2040.6
$ws.Range("L88").Value = 2040.6
$ws.Range("N88").Value = -2852.6
$ws.Range("H89").Value = 8394.643
$ws.Range("I89").Value = 9010.875
$ws.Range("J89").Value = 7573
$ws.Range("K89").Value = 45054.375
$ws.Range("L89").Value = 37865
$ws.Range("M89").Value = -39438.375
$ws.Range("N89").Value = -49097
$ws.Range("H91").Value = 3800.3333
$ws.Range("J91").Value = 2040.6
$ws.Range("L91").Value = 2040.6
$ws.Range("N91").Value = -4848.6
$ws.Range("H92").Value = 1595.9445
$ws.Range("I92").Value = 2211.6365
$ws.Range("K92").Value = 2211.6365
$ws.Range("M92").Value = -963.6365000000001
$ws.Range("H127").Value = 204600
$ws.Range("I127").Value = 204600
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 613800
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -608840
$ws.Range("N127").ClearContents()
$ws.Range("H131").Value = 2153.4
$ws.Range("I131").Value = 2153.4
$ws.Range("K131").Value = 6460.200000000001
$ws.Range("M131").Value = -1420.200000000001
$ws.Range("H137").Value = 18524334
$ws.Range("I137").Value = 19232962
$ws.Range("K137").Value = 57698886
$ws.Range("M137").Value = -57696336
$ws.Range("H138").Value = 8403.92
$ws.Range("I138").Value = 8085.5713
$ws.Range("J138").Value = 8527.723
$ws.Range("K138").Value = 24256.7139
$ws.Range("L138").Value = 25583.169
$ws.Range("M138").Value = -19116.7139
$ws.Range("N138").Value = -35863.169

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1015836.1
$ws.Range("I32").Value = 1140090.4
$ws.Range("J32").Value = 21802.4
$ws.Range("K32").Value = 1140090.4
$ws.Range("L32").Value = 21802.4
$ws.Range("M32").Value = -1139803.4
$ws.Range("N32").Value = -22376.4
$ws.Range("H45").Value = 1642.1
$ws.Range("I45").Value = 1174.5714
$ws.Range("K45").Value = 1174.5714
$ws.Range("M45").Value = -797.5714
$ws.Range("H61").Value = 3705796
$ws.Range("I61").Value = 2195.68
$ws.Range("J61").Value = 50000800
$ws.Range("K61").Value = 2195.68
$ws.Range("L61").Value = 50000800
$ws.Range("M61").Value = -1983.68
$ws.Range("N61").Value = -50001224
$ws.Range("H74").Value = 1067845.8
$ws.Range("I74").Value = 1211713.4
$ws.Range("K74").Value = 1211713.4
$ws.Range("M74").Value = -1210839.4
$ws.Range("H77").Value = 1067845.8
$ws.Range("I77").Value = 1211713.4
$ws.Range("K77").Value = 6058567
$ws.Range("M77").Value = -6054199
$ws.Range("H88").Value = 1707.7
$ws.Range("J88").Value = 2244.8667
$ws.Range("L88").Value = 2244.8667
$ws.Range("N88").Value = -3056.8667
$ws.Range("H91").Value = 1707.7
$ws.Range("J91").Value = 2244.8667
$ws.Range("L91").Value = 2244.8667
$ws.Range("N91").Value = -5052.8667
$ws.Range("H102").Value = 4179.9375
$ws.Range("I102").Value = 4192
$ws.Range("K102").Value = 4192
$ws.Range("M102").Value = -2570
$ws.Range("H136").Value = 3705796
$ws.Range("I136").Value = 2195.68
$ws.Range("J136").Value = 50000800
$ws.Range("K136").Value = 6587.039999999999
$ws.Range("L136").Value = 150002400
$ws.Range("M136").Value = -4037.039999999999
$ws.Range("N136").Value = -150007500

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1517.52
$ws.Range("I107").Value = 1180.8182
$ws.Range("J107").Value = 1782.0714
$ws.Range("K107").Value = 1180.8182
$ws.Range("L107").Value = 1782.0714
$ws.Range("M107").Value = 739.1818000000001
$ws.Range("N107").Value = -5622.0714

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1806668.2
$ws.Range("I31").Value = 2083617.1
$ws.Range("K31").Value = 2083617.1
$ws.Range("M31").Value = -2083322.1
$ws.Range("H34").Value = 1806668.2
$ws.Range("I34").Value = 2083617.1
$ws.Range("K34").Value = 2083617.1
$ws.Range("M34").Value = -2083415.1
$ws.Range("H41").Value = 21139.75
$ws.Range("J41").Value = 23350
$ws.Range("L41").Value = 23350
$ws.Range("N41").Value = -24206
$ws.Range("H86").Value = 9715.954
$ws.Range("I86").Value = 6250
$ws.Range("J86").Value = 10062.55
$ws.Range("K86").Value = 6250
$ws.Range("L86").Value = 10062.55
$ws.Range("M86").Value = -5127
$ws.Range("N86").Value = -12308.55
$ws.Range("H89").Value = 9715.954
$ws.Range("I89").Value = 6250
$ws.Range("J89").Value = 10062.55
$ws.Range("K89").Value = 31250
$ws.Range("L89").Value = 50312.75
$ws.Range("M89").Value = -25634
$ws.Range("N89").Value = -61544.75
$ws.Range("H107").Value = 454.8
$ws.Range("I107").Value = 396.29413
$ws.Range("K107").Value = 396.29413
$ws.Range("M107").Value = 1523.70587
$ws.Range("H132").Value = 2083.738
$ws.Range("I132").Value = 1947.9445
$ws.Range("J132").Value = 2898.5
$ws.Range("K132").Value = 5843.833500000001
$ws.Range("L132").Value = 8695.5
$ws.Range("M132").Value = -3313.833500000001
$ws.Range("N132").Value = -13755.5
$ws.Range("H134").Value = 3694.0417
$ws.Range("I134").Value = 1901.8684
$ws.Range("K134").Value = 5705.6052
$ws.Range("M134").Value = -3170.6052

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 13129.733
$ws.Range("I120").Value = 8194.6
$ws.Range("K120").Value = 24583.8
$ws.Range("M120").Value = -19745.8
$ws.Range("H132").Value = 1237.3334
$ws.Range("I132").Value = 1271.2858
$ws.Range("J132").Value = 1189.8
$ws.Range("K132").Value = 11441.5722
$ws.Range("L132").Value = 10708.2
$ws.Range("M132").Value = -8911.572200000001
$ws.Range("N132").Value = -15768.2
$ws.Range("H137").Value = 6244.16
$ws.Range("J137").Value = 10114.77
$ws.Range("L137").Value = 30344.31
$ws.Range("N137").Value = -40544.31

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1012
$ws.Range("I102").Value = 832.36365
$ws.Range("K102").Value = 832.36365
$ws.Range("M102").Value = 789.63635

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 33265
$ws.Range("J47").Value = 33265
$ws.Range("L47").Value = 33265
$ws.Range("N47").Value = -34245
$ws.Range("H52").Value = 33265
$ws.Range("J52").Value = 33265
$ws.Range("L52").Value = 33265
$ws.Range("N52").Value = -33731
$ws.Range("H55").Value = 1428.6052
$ws.Range("I55").Value = 1210.8823
$ws.Range("K55").Value = 1210.8823
$ws.Range("M55").Value = -1037.8823
$ws.Range("H61").Value = 17444.223
$ws.Range("I61").Value = 16333.167
$ws.Range("J61").Value = 19666.334
$ws.Range("K61").Value = 16333.167
$ws.Range("L61").Value = 19666.334
$ws.Range("M61").Value = -16131.167
$ws.Range("N61").Value = -20070.334
$ws.Range("H113").Value = 17444.223
$ws.Range("I113").Value = 16333.167
$ws.Range("J113").Value = 19666.334
$ws.Range("K113").Value = 16333.167
$ws.Range("L113").Value = 19666.334
$ws.Range("M113").Value = -14163.167
$ws.Range("N113").Value = -24006.334
$ws.Range("H132").Value = 1518353
$ws.Range("I132").Value = 6668846.5
$ws.Range("J132").Value = 3502
$ws.Range("K132").Value = 20006539.5
$ws.Range("L132").Value = 10506
$ws.Range("M132").Value = -20004009.5
$ws.Range("N132").Value = -15566

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 8482.5
$ws.Range("J53").Value = 8482.5
$ws.Range("L53").Value = 8482.5
$ws.Range("N53").Value = -9696.5
$ws.Range("H122").Value = 23533.51
$ws.Range("J122").Value = 142590.62
$ws.Range("L122").Value = 427771.86
$ws.Range("N122").Value = -432671.86
$ws.Range("H132").Value = 2978329
$ws.Range("I132").Value = 3207032.5
$ws.Range("J132").Value = 5181.5
$ws.Range("K132").Value = 9621097.5
$ws.Range("L132").Value = 15544.5
$ws.Range("M132").Value = -9618567.5
$ws.Range("N132").Value = -20604.5
